$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 292. Excel copies formatting from the row
# above (row 291) by default when inserting, matching the style (s="2")
# already used by the neighboring date cells in column D.
$ws.Rows.Item(292).Insert()

# The newly inserted row 292 starts as a duplicate data record (same
# market/product) with an updated date and updated price figures.
$ws.Range("A292").Value2 = 3
$ws.Range("B292").Value2 = "Femacal de La Calera"
$ws.Range("C292").Value2 = "Coquimbo"
$ws.Range("D292").Value2 = 44694
$ws.Range("E292").Value2 = 5
$ws.Range("F292").Value2 = 100112009
$ws.Range("G292").Value2 = "Acelga"
$ws.Range("H292").Value2 = "Sin especificar"
$ws.Range("I292").Value2 = "Primera"
$ws.Range("J292").Value2 = 250
$ws.Range("K292").Value2 = 2500
$ws.Range("L292").Value2 = 3000
$ws.Range("M292").Value2 = 2760
$ws.Range("N292").Value2 = "$/docena de atados (6 kilos)"
$ws.Range("O292").Value2 = "Provincia de Quillota"
$ws.Range("P292").Value2 = 460
$ws.Range("Q292").Value2 = 6
$ws.Range("R292").Value2 = "Hortaliza"
